# Se automatiza el flujo cobro recibo
#
# - Adds a new worksheet "Hoja2" right after "Hoja1".
# - Populates row 1 of "Hoja2" with the new category breadcrumb strings.
# - Makes "Hoja2" the active/selected sheet (tabSelected + activeTab).
# - Leaves a selection on "Hoja1" at C2 (no longer the selected tab).
# - Selects C1 on "Hoja2".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet immediately after "Hoja1".
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Hoja2"

$ws2.Range("A1").Value = "Acsel;Area Administrativa y Financiera"
$ws2.Range("B1").Value = "Acsel;Area Administrativa y Financiera;Cobros - Ingresos"
$ws2.Range("C1").Value = "Acsel;Area Administrativa y Financiera;Cobros - Ingresos;Ingresos"

# Restore the non-active selection left on Hoja1 (C2).
$ws1.Range("C2").Select() | Out-Null

# Hoja2 becomes the active sheet/tab, selection resting on C1.
$ws2.Activate() | Out-Null
$ws2.Range("C1").Select() | Out-Null
